$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Seed new shared-string table entries in the exact order they were
#      originally authored, via a scratch sheet that is removed afterwards ----
$scratch = $wb.Worksheets.Add()
$scratch.Range("A1").Value = "Edit a Board for a Project"
$scratch.Range("A1").Value = "View Active Tasks"
$scratch.Range("A1").Value = "View Available Tasks"
$scratch.Range("A1").Value = "Update Stage of an active Task"
$scratch.Range("A1").Value = "Update Stage of my assigned Task"
$scratch.Range("A1").Value = "Post a comment to an active Task"
$scratch.Range("A1").Value = "Create/Edit a Task on the Board of a Project for my Group"
$scratch.Delete()

# ---- Insert 4 additional rows (20 data rows -> 24 data rows) ----
$ws.Rows("19:22").Insert()

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Admin"
$ws.Range("C2").Value = "View, Create, and Delete Employee"
$ws.Range("D2").Value = "I can manage the set of employees in the system"

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "Admin"
$ws.Range("C3").Value = "View, Create, and Delete Groups"
$ws.Range("D3").Value = "I can organize employees into collaborative groups with a manager supervising work for the group"

# Row 4
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "Manager*"
$ws.Range("C4").Value = "View, Create, and Delete Projects"
$ws.Range("D4").Value = "I can assign work to groups"

# Row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "Manager*"
$ws.Range("C5").Value = "Edit Title and/or Description of Project"
$ws.Range("D5").Value = "I can assign a name and description for the project to help others understand the purpose of the project"

# Row 6
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = "Manager*"
$ws.Range("C6").Value = "Edit a Board for a Project"
$ws.Range("D6").Value = "I can define the Stages of work expected for a project and setup for any handoffs needed during project work"

# Row 7
$ws.Range("A7").Value = 0.5
$ws.Range("B7").Value = "Admin"
$ws.Range("C7").Value = "Promote/Demote Employee to/from Manager for a Group"
$ws.Range("D7").Value = "I can manage who is supervising a particular group of employees"

# Row 8
$ws.Range("A8").Value = 0.5
$ws.Range("B8").Value = "Manager*"
$ws.Range("C8").Value = "Assign/Remove Group to/from Project"
$ws.Range("D8").Value = "I can control who is actively collaborating to complete a project"

# Row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Employee**"
$ws.Range("C9").Value = "View Board for Project"
$ws.Range("D9").Value = "I can overview work in progress and overall status for a project"

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Employee**"
$ws.Range("C10").Value = "Create/Edit a Task on the Board of a Project for my Group"
$ws.Range("D10").Value = "I can define work needed for a Project for my Group"

# Row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Employee**"
$ws.Range("C11").Value = "Update Stage of my assigned Task"
$ws.Range("D11").Value = "I can notify others about the progress on a Task"

# Row 12
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Employee**"
$ws.Range("C12").Value = "Assign a Task on a Board from a Project for my Group to an Employee in my Group"
$ws.Range("D12").Value = "I can delegate work to an Employee in my Group"

# Row 13
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Employee**"
$ws.Range("C13").Value = "View Active Tasks"
$ws.Range("D13").ClearContents()
$ws.Range("A13:D13").Interior.Color = 65535

# Row 14
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Employee**"
$ws.Range("C14").Value = "View Available Tasks"
$ws.Range("D14").ClearContents()
$ws.Range("A14:D14").Interior.Color = 65535

# Row 15
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "Employee**"
$ws.Range("C15").Value = "Update Stage of an active Task"
$ws.Range("D15").ClearContents()
$ws.Range("A15:D15").Interior.Color = 65535

# Row 16
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Employee**"
$ws.Range("C16").Value = "Post a comment to a Task on the Board of a Project for my Group"
$ws.Range("D16").Value = "I can discuss the work for a Task in a space that will be preserved for later review"

# Row 17
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Employee**"
$ws.Range("C17").Value = "Post a comment to an active Task"
$ws.Range("D17").Value = "I can discuss the work for a Task in a space that will be preserved for later review"
$ws.Range("A17:D17").Interior.Color = 65535

# Row 18
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Manager*"
$ws.Range("C18").Value = "View Task Status Report for Employee in my Group"
$ws.Range("D18").Value = "I can monitor the progress for a single Employee"

# Row 19
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "Manager*"
$ws.Range("C19").Value = "View Task Status Report for a Project I Manage"
$ws.Range("D19").Value = "I can monitor the progress for a Project"

# Row 20
$ws.Range("A20").Value = 3
$ws.Range("B20").Value = "Manager*"
$ws.Range("C20").Value = "View Task Status Report for a Group I manage"
$ws.Range("D20").Value = "I can monitor the progress for a single Group"

# Row 21
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Manager*"
$ws.Range("C21").Value = "Assign/Remove Group(s) to/from Stage of Project"
$ws.Range("D21").Value = "I can define who will work on the project in a particular stage"

# Row 22
$ws.Range("A22").Value = 4
$ws.Range("B22").Value = "Manager*"
$ws.Range("C22").Value = "Assign an Employee in my Group as a reporter for a Task on a Board for my Group"
$ws.Range("D22").Value = "I can define who will be notified about updates for a task"

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "*"
$ws.Range("C23").Value = "All features available to Managers are also available to Admins"
$ws.Range("D23").ClearContents()

# Row 24
$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "**"
$ws.Range("C24").Value = "All features available to Employees are also available to Managers and Admins"
$ws.Range("D24").ClearContents()

# ---- Sort state metadata (range grew from D18 to D27) ----
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("A2:A27"))
$so.SortFields.Add($ws.Range("B2:B27"))
$so.SetRange($ws.Range("A2:D27"))
$so.Header = 0
$so.Apply()

# ---- View: zoom + selection ----
$ws.Application.ActiveWindow.Zoom = 80
$ws.Range("F12").Select()

